$d = $word.ActiveDocument

# The document has one real footnote and one real endnote (the
# separator / continuation-separator notes are not part of these
# collections). Both currently start their paragraph with:
#   <w:tabs><w:tab w:val="start" w:pos="400"/></w:tabs>
# and need that replaced with:
#   <w:ind w:start="400" w:hanging="400"/>
#
# Setting ParagraphFormat.LeftIndent/FirstLineIndent individually only
# ever produces a single attribute on <w:ind> (each assignment rewrites
# the whole element), so instead we rebuild the note's paragraph via
# Range.InsertXML with the full desired WordprocessingML, which lets us
# specify both w:start and w:hanging together on one <w:ind/>.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$footnote = $d.Footnotes.Item(1)
$footnoteXml = '<w:p ' + $wNs + '><w:pPr><w:ind w:start="400" w:hanging="400"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:footnoteRef/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">This is the footnote content.</w:t></w:r></w:p>'
$footnote.Range.InsertXML($footnoteXml)

$endnote = $d.Endnotes.Item(1)
$endnoteXml = '<w:p ' + $wNs + '><w:pPr><w:ind w:start="400" w:hanging="400"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:endnoteRef/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">This is the endnote content.</w:t></w:r></w:p>'
$endnote.Range.InsertXML($endnoteXml)
